$d = $word.ActiveDocument

# Locate the title cell ("CU14.2 - Eliminar un Grupo") in the first table.
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(1, 2)
$cellRange = $cell.Range
$cellText = $cellRange.Text
$grupoOffset = $cellText.IndexOf("Grupo")

$grupoStart = $cellRange.Start + $grupoOffset
$grupoEnd = $grupoStart + 5          # length of "Grupo"
$membresiaEnd = $grupoStart + 9      # length of "Membresía"
$unStart = $grupoStart - 3           # length of "un " before "Grupo"

# 1) Replace "Grupo" with "Membresía" (stays inside the existing run for now).
$grupoRange = $d.Range($grupoStart, $grupoEnd)
$grupoRange.Text = "Membresía"

# 2) Re-establish the run boundary between the preceding " " run and "un
#    Membresía" by nudging (and restoring) a character property - this
#    forces the engine to split the text back into separate runs instead
#    of leaving everything coalesced into the space-run.
$wideRange = $d.Range($unStart, $membresiaEnd)
$wideRange.Font.Bold = $false
$wideRange.Font.Bold = $true

# 3) Split "Membresía" away from "un " the same way, so the final text
#    reads "un " + "Membresía" as two separate (still bold) runs.
$membresiaRange = $d.Range($grupoStart, $membresiaEnd)
$membresiaRange.Font.Bold = $false
$membresiaRange.Font.Bold = $true
